$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Travail effectué" text for Jeudi (row 6) to append ". Carte"
$ws.Range("F6").Value = "Modifications pour relance facture et pdf relance facture. Carte"

# Add the new "Pause" entry for Jeudi (row 6)
$ws.Range("C6").Value = "12:16 - 13:45 / 18:15 - 19:25"
$ws.Range("C6").NumberFormat = $ws.Range("B6").NumberFormat

# Update the active selection to match
$ws.Range("C7").Select()
